# Update cryptos list values (price & 1h volume change) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells use "General" number formatting. Force each cell to text format
# first so numeric-looking strings (e.g. "610.50", "24.00") keep their original
# text representation instead of being auto-coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.673.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.621.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.619.25"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.417"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.234.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.614.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.757.91"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.34"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.618.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.160"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.45"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.90"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.64"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.41"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0864"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.70%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.969"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.70%  "
